# Regenerate the handoff/handback status report timestamps.
#
# Two "Latest Handback DateTime" / "Latest Handoff Date" timestamps were
# refreshed by a re-run of the report generator:
#   2016-03-22 10:28:40 -> 2016-03-22 10:29:35  (Overview + de-de sheets)
#   2016-03-22 10:28:36 -> 2016-03-22 10:29:31  (zh-cn sheet)
#
# Two files (36e8d332-... and c5198dc3-...) that previously had their own
# distinct timestamp (2016-03-22 10:29:09 / 2016-03-22 10:29:05) now land
# on the very same refreshed timestamp as aea18086-...md (row 7), so all
# of those cells converge on the single updated value per sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets("Overview")
$zhcn     = $wb.Worksheets("zh-cn")
$dede     = $wb.Worksheets("de-de")

# Overview sheet: column D ("Latest Handoff Date") rows 7,10-16
$overview.Range("D7").Value  = "2016-03-22 10:29:35"
$overview.Range("D10").Value = "2016-03-22 10:29:35"
$overview.Range("D11").Value = "2016-03-22 10:29:35"
$overview.Range("D12").Value = "2016-03-22 10:29:35"
$overview.Range("D13").Value = "2016-03-22 10:29:35"
$overview.Range("D14").Value = "2016-03-22 10:29:35"
$overview.Range("D15").Value = "2016-03-22 10:29:35"
$overview.Range("D16").Value = "2016-03-22 10:29:35"

# zh-cn sheet: column E ("Latest Handoff Datetime") rows 7,10-16
$zhcn.Range("E7").Value  = "2016-03-22 10:29:31"
$zhcn.Range("E10").Value = "2016-03-22 10:29:31"
$zhcn.Range("E11").Value = "2016-03-22 10:29:31"
$zhcn.Range("E12").Value = "2016-03-22 10:29:31"
$zhcn.Range("E13").Value = "2016-03-22 10:29:31"
$zhcn.Range("E14").Value = "2016-03-22 10:29:31"
$zhcn.Range("E15").Value = "2016-03-22 10:29:31"
$zhcn.Range("E16").Value = "2016-03-22 10:29:31"

# de-de sheet: column E ("Latest Handoff Datetime") rows 7,10-16
$dede.Range("E7").Value  = "2016-03-22 10:29:35"
$dede.Range("E10").Value = "2016-03-22 10:29:35"
$dede.Range("E11").Value = "2016-03-22 10:29:35"
$dede.Range("E12").Value = "2016-03-22 10:29:35"
$dede.Range("E13").Value = "2016-03-22 10:29:35"
$dede.Range("E14").Value = "2016-03-22 10:29:35"
$dede.Range("E15").Value = "2016-03-22 10:29:35"
$dede.Range("E16").Value = "2016-03-22 10:29:35"
